$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62: a few fields were corrected in place -------------------------
$ws.Range("A62").Value = 80976130
$ws.Range("L62").Value = "hane"
$ws.Range("M62").Value = "frispringande/krypande"
$ws.Range("AC62").Value = "Ny lokal, och blott tredje kända i Södermanland!? Närmast funnen i Tyresta NP. Bör eftersökas på fler lokaler i kommunen! Grävde först fram fragment av en ad hona, men kunde sedan finna en vuxen hane i en perfekt rödmurken granlåga i sent nedbrytningsstadium."

# --- Rows 63-66: the four observation records got reordered ---------------
# New row 63 <- old row 64 (Rödgul trumpetsvamp)
$ws.Range("A63").Value = 111683845
$ws.Range("B63").Value = 89183
$ws.Range("E63").Value = 3215
$ws.Range("F63").Value = "Rödgul trumpetsvamp"
$ws.Range("G63").Value = "Craterellus lutescens"
$ws.Range("H63").Value = "(Fr.) Fr."
$ws.Range("I63").Value = ""
$ws.Range("J63").Value = ""
$ws.Range("P63").Value = "Fiskarsundet, Srm"
$ws.Range("Q63").Value = 689111.5690902721
$ws.Range("R63").Value = 6570305.953062683
$ws.Range("S63").Value = 23
$ws.Range("Z63").Value = "09:36"
$ws.Range("AB63").Value = "09:36"

# New row 64 <- old row 66 (Fjällig taggsvamp s.str.)
$ws.Range("A64").Value = 111683853
$ws.Range("B64").Value = 90687
$ws.Range("E64").Value = 5964
$ws.Range("F64").Value = "Fjällig taggsvamp s.str."
$ws.Range("G64").Value = "Sarcodon imbricatus s.str."
$ws.Range("H64").Value = "(L.:Fr.) P.Karst."
$ws.Range("Z64").Value = "09:34"
$ws.Range("AB64").Value = "09:34"

# New row 65 <- old row 63 (Svavelriska)
$ws.Range("A65").Value = 111683850
$ws.Range("B65").Value = 90332
$ws.Range("E65").Value = 4769
$ws.Range("F65").Value = "Svavelriska"
$ws.Range("G65").Value = "Lactarius scrobiculatus"
$ws.Range("H65").Value = "(Scop.:Fr.) Fr."
$ws.Range("I65").NumberFormat = "@"
$ws.Range("I65").Value = "3"
$ws.Range("J65").Value = "fruktkroppar"
$ws.Range("L65").ClearContents()
$ws.Range("P65").Value = "Bergaholm, Tyresö kn, Srm"
$ws.Range("Q65").Value = 689075.4602011892
$ws.Range("R65").Value = 6570319.534944151
$ws.Range("S65").Value = 20
$ws.Range("Z65").Value = "09:25"
$ws.Range("AB65").Value = "09:25"

# New row 66 <- old row 65 (Sårläka)
$ws.Range("A66").Value = 111683856
$ws.Range("B66").Value = 108219
$ws.Range("E66").Value = 219711
$ws.Range("F66").Value = "Sårläka"
$ws.Range("G66").Value = "Sanicula europaea"
$ws.Range("H66").Value = "L."
$ws.Range("L66").Value = ""
